# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 1608
    5  = 616
    6  = 1095
    7  = 14
    8  = 11442
    9  = 27
    10 = 91
    15 = 12366
    16 = 13039
    18 = 139
    21 = 18
    24 = 105
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
